$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.771.03'
$ws.Range("E2").Value = '  -2.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.958.78'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.08'
$ws.Range("E5").Value = '  -2.62%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4760'
$ws.Range("E7").Value = '  -5.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4023'
$ws.Range("E8").Value = '  -5.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.80'
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08471'
$ws.Range("E10").Value = '  -5.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.059'
$ws.Range("E11").Value = '  -5.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.35'
$ws.Range("E12").Value = '  -4.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.963.61'
$ws.Range("E13").Value = '  -4.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.606'
$ws.Range("E14").Value = '  -6.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.214'
$ws.Range("E15").Value = '  -4.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.013'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.08'
$ws.Range("E17").Value = '  -5.85%  '
$ws.Range("E18").Value = '  -4.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06601'
$ws.Range("E19").Value = '  -1.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.52'
$ws.Range("E20").Value = '  -6.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.791'
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.784.93'
$ws.Range("E23").Value = '  -2.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.50'
$ws.Range("E24").Value = '  -4.45%  '
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.235.37'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.15'
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.14'
$ws.Range("E28").Value = '  -3.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.909'
$ws.Range("E29").Value = '  -7.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.148'
$ws.Range("E30").Value = '  -7.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '123.70'
$ws.Range("E31").Value = '  -4.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.001'
$ws.Range("E32").Value = '  -5.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09567'
$ws.Range("E33").Value = '  -4.05%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.658'
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.445'
$ws.Range("E35").Value = '  -8.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.671'
$ws.Range("E36").Value = '  -3.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02343'
$ws.Range("E37").Value = '  -5.47%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06207'
$ws.Range("E38").Value = '  -2.58%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.262'
$ws.Range("E39").Value = '  -4.25%  '
$ws.Range("E40").Value = '  -7.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6205'
$ws.Range("E41").Value = '  -5.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.07'
$ws.Range("E42").Value = '  -5.60%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1912'
$ws.Range("E44").Value = '  -7.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.331'
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5935'
$ws.Range("E46").Value = '  -6.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.89'
$ws.Range("E47").Value = '  -5.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.064'
$ws.Range("E48").Value = '  -6.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.418'
$ws.Range("E49").Value = '  -3.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000334'
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06825'
$ws.Range("E51").Value = '  -2.53%  '
